# 🚌 141: 31/12 18:59 LP1912+6203+6173
# Appends newly-scraped arrival rows to the "LP1912" and "6203-6173" sheets
# and refreshes the "Última actualización" / "Total filas" header cells on
# all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": 17 new rows (1055-1071), columns B:G
#   B=Hora_Scrap  C=Hora_Llegada  D=Línea  E=Minutos  F=Parada  G=Fecha
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 15:59:25"
$ws1.Range("A3").Value = "Total filas: 1070"

$data1 = @(
    @("15:59:15", "16:04", "23_HERNANDEZ",         5, "LP1912", "31/12/2025"),
    @("15:59:15", "16:10", "16_SANTA ANA",        11, "LP1912", "31/12/2025"),
    @("15:59:15", "16:15", "10_OLMOS",            16, "LP1912", "31/12/2025"),
    @("15:59:15", "16:22", "16_SANTA ANA",        23, "LP1912", "31/12/2025"),
    @("15:59:15", "16:25", "11_ETCHEVERRY",       26, "LP1912", "31/12/2025"),
    @("15:59:15", "16:32", "16_P MOR-SANTA ANA",  33, "LP1912", "31/12/2025"),
    @("15:59:15", "16:34", "23_HERNANDEZ",        35, "LP1912", "31/12/2025"),
    @("15:59:15", "16:37", "17X38_ROMERO",        38, "LP1912", "31/12/2025"),
    @("15:59:15", "16:52", "16_SANTA ANA",        53, "LP1912", "31/12/2025"),
    @("15:59:15", "16:54", "10_OLMOS",            55, "LP1912", "31/12/2025"),
    @("15:59:15", "17:05", "14_ABASTO",           66, "LP1912", "31/12/2025"),
    @("15:59:15", "17:08", "15_ABASTO",           69, "LP1912", "31/12/2025"),
    @("15:59:15", "17:25", "11_ETCHEVERRY",       86, "LP1912", "31/12/2025"),
    @("15:59:15", "17:28", "15_ABASTO",           89, "LP1912", "31/12/2025"),
    @("15:59:15", "17:34", "23_HERNANDEZ",        95, "LP1912", "31/12/2025"),
    @("15:59:15", "17:35", "10_OLMOS",            96, "LP1912", "31/12/2025"),
    @("15:59:15", "17:36", "16_P MOR-SANTA ANA",  97, "LP1912", "31/12/2025")
)

$startRow1 = 1055
$rowCount1 = $data1.Count
$colCount1 = 6
$arr1 = New-Object 'object[,]' $rowCount1, $colCount1
for ($i = 0; $i -lt $rowCount1; $i++) {
    for ($j = 0; $j -lt $colCount1; $j++) {
        $arr1[$i, $j] = $data1[$i][$j]
    }
}
$endRow1 = $startRow1 + $rowCount1 - 1
$ws1.Range("B$startRow1`:G$endRow1").Value = $arr1

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": no new rows, just the shared timestamp refresh
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 31/12/2025 15:59:25"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": 2 new rows (132-133), columns B:G
#   B=Fecha  C=Hora_Scrap  D=Hora_Llegada  E=Línea  F=Minutos  G=Parada
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 15:59:25"
$ws3.Range("A3").Value = "Total filas: 132"

$data3 = @(
    @("31/12/2025", "15:59:25", "16:26", "215B_LP-P MOR-1 Y 57", 27, "L6173"),
    @("31/12/2025", "15:59:20", "16:58", "215C_LA PLATA",        59, "L6203")
)

$startRow3 = 132
$rowCount3 = $data3.Count
$colCount3 = 6
$arr3 = New-Object 'object[,]' $rowCount3, $colCount3
for ($i = 0; $i -lt $rowCount3; $i++) {
    for ($j = 0; $j -lt $colCount3; $j++) {
        $arr3[$i, $j] = $data3[$i][$j]
    }
}
$endRow3 = $startRow3 + $rowCount3 - 1
$ws3.Range("B$startRow3`:G$endRow3").Value = $arr3

Write-Host "Applied 141 scrape update: +$rowCount1 rows to LP1912, +$rowCount3 rows to 6203-6173"
